$wb = $excel.ActiveWorkbook

# --- Update the "Date" value on the Metadata sheet ---
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2023-07-27T13:09:10+00:00"

# --- Update the "prolongee-prorogee" concept to "prorogee" on the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("B5").Value = "prorogee"
$concepts.Range("C5").Value = "Prorogée"
